$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fileName = "PVP_Baocaothuongnien_2022.pdf"
$keywords = @("Công nghệ", "Đảm bảo", "Công việc", "Làm việc", "AI", "Chuyển đổi số", "Tra cứu", "Blockchain")

$startRow = 18
for ($i = 0; $i -lt $keywords.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $fileName
    $ws.Cells.Item($row, 2).Value = $keywords[$i]
    $ws.Cells.Item($row, 3).Value = 0
}
